$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Michael Beaver")

# --- Row 27: copy formatting from row 24 (same target styles) ---
$ws.Range("A24:I24").Copy() | Out-Null
$ws.Range("A27:I27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A27").Value = 41683
$ws.Range("B27").Value = "Updated the software requirements specification document as per Travis's comments and suggestions. Drafted programming language proposal."
$ws.Range("I27").Value = 0.75

# --- Row 28: copy formatting from row 7 (same target styles) ---
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A28:I28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A28").Value = 41683
$ws.Range("B28").Value = "Client meeting. Recorded meeting minutes and client's responses to questions. Transcribed team members' reports of client's answers onto the meeting minutes."
$ws.Range("I28").Value = 1.25

$excel.CutCopyMode = 0

# --- Update view / selection to match final state ---
$ws.Range("B30:H30").Select() | Out-Null

$wb.Save()
